# Atualizado por script em 05-11-2023 14:45
#
# This script:
#  1) Rewrites the odds/score data (columns F..V) for the matches whose
#     betting odds were re-scraped (rows 323-331 and 338-339), while the
#     match index (A), pais/torneio/temporada (B-D) and kickoff date (E)
#     columns are left untouched.
#  2) Appends a brand-new match row (341): FC Osaka vs Fukushima United.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row-by-row replacement data: F, G, H, I, J, L, M, N, P, Q, R, T, U, V
#    (columns A-E and K/O/S stay exactly as they already are).
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=323; F="Osaka";               G=1; H="Gainare Tottori";      I=1; J=2.26; L=2.36; M="28/10/2023 06:36"; N=3.09; P=3.18; Q="28/10/2023 06:58"; R=2.99; T=3.14; U="28/10/2023 06:58"; V="https://www.betexplorer.com/football/japan/j3-league/fc-osaka-gainare-tottori/zmCDkxpH/" },
    @{ Row=324; F="Sagamihara";          G=0; H="Yamaga";               I=2; J=2.77; L=3.67; M="28/10/2023 06:16"; N=3.08; P=3.42; Q="28/10/2023 06:16"; R=2.41; T=2.03; U="28/10/2023 06:16"; V="https://www.betexplorer.com/football/japan/j3-league/sc-sagamihara-matsumoto-yamaga/ddD9jIUA/" },
    @{ Row=325; F="Imabari";             G=0; H="Grulla Morioka";       I=1; J=2.06; L=1.96; M="29/10/2023 03:23"; N=3.59; P=3.39; Q="29/10/2023 03:23"; R=3.06; T=3.94; U="29/10/2023 03:23"; V="https://www.betexplorer.com/football/japan/j3-league/imabari-iwate-grulla-morioka/6uMs6Gic/" },
    @{ Row=326; F="Tegevajaro Miyazaki"; G=0; H="Ryukyu";               I=2; J=2.31; L=2.71; M="29/10/2023 04:40"; N=3.1;  P=3.31; Q="29/10/2023 03:10"; R=2.9;  T=2.59; U="29/10/2023 04:40"; V="https://www.betexplorer.com/football/japan/j3-league/tegevajaro-miyazaki-ryukyu/E5Pk4fM9/" },
    @{ Row=327; F="Fukushima Utd";       G=0; H="Kagoshima Utd";        I=1; J=3.17; L=3.62; M="29/10/2023 04:56"; N=3.21; P=3.8;  Q="29/10/2023 04:56"; R=2.11; T=1.93; U="29/10/2023 04:56"; V="https://www.betexplorer.com/football/japan/j3-league/fukushima-united-kagoshima-united/lfhlqzNp/" },
    @{ Row=328; F="Giravanz Kitakyushu"; G=0; H="YSCC";                 I=2; J=2.47; L=2.89; M="29/10/2023 04:12"; N=3.05; P=3.15; Q="29/10/2023 04:12"; R=2.72; T=2.54; U="29/10/2023 04:12"; V="https://www.betexplorer.com/football/japan/j3-league/giravanz-kitakyushu-yscc-yokohama/02Lo5z73/" },
    @{ Row=329; F="Nagano";              G=3; H="Azul Claro Numazu";    I=0; J=2.66; L=2.69; M="29/10/2023 04:29"; N=3.07; P=3.62; Q="29/10/2023 05:58"; R=2.51; T=2.45; U="29/10/2023 04:29"; V="https://www.betexplorer.com/football/japan/j3-league/nagano-parceiro-azul-claro-numazu/Spihrfxi/" },
    @{ Row=330; F="Toyama";              G=1; H="Gifu";                 I=1; J=2.31; L=2.48; M="29/10/2023 05:56"; N=3.08; P=3.29; Q="29/10/2023 05:14"; R=2.91; T=2.85; U="29/10/2023 05:56"; V="https://www.betexplorer.com/football/japan/j3-league/toyama-gifu/SlNw7dxi/" },
    @{ Row=331; F="Vanraure";            G=2; H="Kamatamare";           I=2; J=2.04; L=2.4;  M="29/10/2023 05:58"; N=3.14; P=3.08; Q="29/10/2023 05:58"; R=3.43; T=3.16; U="29/10/2023 05:58"; V="https://www.betexplorer.com/football/japan/j3-league/vanraure-kamatamare-sanuki/WE5MmGET/" },
    @{ Row=338; F="Azul Claro Numazu";   G=0; H="Tegevajaro Miyazaki";  I=1; J=1.63; L=1.89; M="05/11/2023 04:03"; N=3.63; P=3.74; Q="05/11/2023 04:03"; R=4.9;  T=3.82; U="05/11/2023 04:03"; V="https://www.betexplorer.com/football/japan/j3-league/azul-claro-numazu-tegevajaro-miyazaki/lSgZeWTk/" },
    @{ Row=339; F="Kamatamare";          G=2; H="Nagano";               I=3; J=1.88; L=2.71; M="05/11/2023 05:54"; N=3.35; P=3.09; Q="05/11/2023 04:02"; R=3.79; T=2.74; U="05/11/2023 05:54"; V="https://www.betexplorer.com/football/japan/j3-league/kamatamare-sanuki-nagano-parceiro/E5oLJZLd/" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 6).Value  = $r.F   # F home
    $ws.Cells.Item($n, 7).Value  = $r.G   # G home_ft_gols
    $ws.Cells.Item($n, 8).Value  = $r.H   # H away
    $ws.Cells.Item($n, 9).Value  = $r.I   # I away_ft_gols
    $ws.Cells.Item($n, 10).Value = $r.J   # J home_opening_odds
    $ws.Cells.Item($n, 12).Value = $r.L   # L home_closing_odds
    $ws.Cells.Item($n, 13).Value = $r.M   # M home_closing_data_hora
    $ws.Cells.Item($n, 14).Value = $r.N   # N draw_closing_odds
    $ws.Cells.Item($n, 16).Value = $r.P   # P away_opening... closing odds
    $ws.Cells.Item($n, 17).Value = $r.Q   # Q away_closing_data_hora
    $ws.Cells.Item($n, 18).Value = $r.R   # R away_closing_odds
    $ws.Cells.Item($n, 20).Value = $r.T   # T draw_closing (2nd)
    $ws.Cells.Item($n, 21).Value = $r.U   # U closing data hora
    $ws.Cells.Item($n, 22).Value = $r.V   # V url
}

# ---------------------------------------------------------------------
# 2) Append new row 341 (FC Osaka 1 x 0 Fukushima United), copying
#    direct cell formatting (bold/border/alignment on A, date format on
#    E) from row 340 so no new style entries are introduced.
# ---------------------------------------------------------------------
$ws.Range("A340").Copy()
$ws.Range("A341").PasteSpecial(-4122)
$ws.Range("E340").Copy()
$ws.Range("E341").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(341, 1).Value  = 340
$ws.Cells.Item(341, 2).Value  = "japan"
$ws.Cells.Item(341, 3).Value  = "j3-league"
$ws.Cells.Item(341, 4).Value  = "2023"
$ws.Cells.Item(341, 5).Value  = 45235.45833333334
$ws.Cells.Item(341, 6).Value  = "Osaka"
$ws.Cells.Item(341, 7).Value  = 1
$ws.Cells.Item(341, 8).Value  = "Fukushima Utd"
$ws.Cells.Item(341, 9).Value  = 0
$ws.Cells.Item(341, 10).Value = 1.87
$ws.Cells.Item(341, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(341, 12).Value = 2.13
$ws.Cells.Item(341, 13).Value = "05/11/2023 09:20"
$ws.Cells.Item(341, 14).Value = 3.25
$ws.Cells.Item(341, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(341, 16).Value = 3.38
$ws.Cells.Item(341, 17).Value = "05/11/2023 09:20"
$ws.Cells.Item(341, 18).Value = 3.85
$ws.Cells.Item(341, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(341, 20).Value = 3.42
$ws.Cells.Item(341, 21).Value = "05/11/2023 09:20"
$ws.Cells.Item(341, 22).Value = "https://www.betexplorer.com/football/japan/j3-league/fc-osaka-fukushima-united/02cwejqd/"
